$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add CID values for the existing two compounds (new column B data)
$ws.Range("B2").Value = 2776299
$ws.Range("B3").Value = 999

# Add a new row (row 4) for Remdesivir
$ws.Range("A4").Value = "Remdesivir"
$ws.Range("B4").Value = 121304016
$ws.Range("C4").Value = "2-ethylbutyl (2S)-2-[[[(2R,3S,4R,5R)-5-(4-aminopyrrolo[2,1-f][1,2,4]triazin-7-yl)-5-cyano-3,4-dihydroxyoxolan-2-yl]methoxy-phenoxyphosphoryl]amino]propanoate`n"

# D4 must stay text ("602.6`n") rather than being auto-converted to a number
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "602.6`n"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "C27H35N6O8P`n"

# Undo the row-height autofit side effect triggered by the embedded newlines
$ws.Rows(4).AutoFit()
